$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 25.46000000000054
$ws.Range("G2").Value = [double]"3.299190531880214e-08"
$ws.Range("H2").Value = [double]"1.995500767475569e-07"
$ws.Range("K2").Value = 6.149627910164218
$ws.Range("L2").Value = "[3.3845955291209733, 8.914660291207463]"
$ws.Range("M2").Value = [double]"1.626694268508366e-05"
$ws.Range("N2").Value = [double]"1.626694268508366e-05"
$ws.Range("O2").Value = -0.9308422677303092
$ws.Range("P2").Value = "[-1.3711055024676178, -0.49057903299300065]"
$ws.Range("Q2").Value = [double]"4.070705522418727e-05"
$ws.Range("R2").Value = [double]"4.070705522418727e-05"
$ws.Range("S2").Value = 10.36195210474986
$ws.Range("T2").Value = "[8.938652980695535, 11.785251228804183]"
$ws.Range("W2").Value = 3.771851851851931
$ws.Range("X2").Value = 1.987867867867908
$ws.Range("Y2").Value = 5.555835835835953

# Row 3 updates
$ws.Range("E3").Value = 23.96000000000031
$ws.Range("G3").Value = [double]"2.446563618363484e-09"
$ws.Range("H3").Value = [double]"6.782955863648965e-08"
$ws.Range("K3").Value = 5.788568173046826
$ws.Range("L3").Value = "[3.4855376917031418, 8.091598654390511]"
$ws.Range("M3").Value = [double]"1.120391559084055e-06"
$ws.Range("N3").Value = [double]"2.24078311816811e-06"
$ws.Range("O3").Value = 2.72334258058935
$ws.Range("P3").Value = "[2.3082372449798885, 3.1384479161988117]"
$ws.Range("S3").Value = 10.36783102588813
$ws.Range("T3").Value = "[9.120039177614995, 11.615622874161264]"
$ws.Range("W3").Value = 13.57493493493511
$ws.Range("X3").Value = 11.99199199199215
$ws.Range("Y3").Value = 15.15787787787807
